$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A11").Value = "22.12.2023"
$ws.Range("D11").Value = "13:00-15:00; 16:30-"
$ws.Range("D11").Select()
